$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update description text for C7 (220uf electrolytic capacitor 15V -> 10V)
$ws.Range("C7").Value = "220uf electrolytic capacitor 10V"

# Update quantities
$ws.Range("D3").Value = 10
$ws.Range("D5").Value = 3

# Add new BOM row 35 (Resistor 9K31 0.1%)
$ws.Range("C35").Value = "Resistor 9K31 0.1%"
$ws.Range("D35").Value = 2
$ws.Range("E35").Value = "RS"
$ws.Range("F35").Value = "754-7095"
$ws.Range("G35").Value = 64
$ws.Range("H35").Value = 128
$ws.Range("I35").Value = "R7, R10"

# Update Total row: move from row 36 formula text + SUM range
$ws.Range("F36").Value = "Total"
$ws.Range("H36").Formula = "=SUM(H2:H35)"

# Update sheet view (scroll so row 5 is the top visible row, select H37)
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H37").Select()
